$wb = $excel.ActiveWorkbook

# --- Update "Sprint 4 Backlog" data on "Sprint 4 Burndown" sheet ---
$ws4 = $wb.Worksheets.Item("Sprint 4 Burndown")

$ws4.Range("C8").Value = 0
$ws4.Range("D8").Value = 0

$ws4.Range("C9").Value = 2
$ws4.Range("D9").Value = 0

$ws4.Range("C10").Value = 1
$ws4.Range("D10").Value = 0

$ws4.Range("C11").Value = 2
$ws4.Range("D11").Value = 2

$ws4.Range("C12").Value = 2
$ws4.Range("D12").Value = 2

# --- Selections / active sheet ---
$ws3 = $wb.Worksheets.Item("Sprint 3 Burndown")
$ws3.Range("D18").Select()

$ws4.Activate()
$ws4.Range("D24").Select()
